$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1524, 1).Value2 = 5
$ws.Cells.Item(1524, 2).Value2 = 2
$ws.Cells.Item(1524, 3).Value2 = 6
$ws.Cells.Item(1524, 4).Value2 = 18

$ws.Cells.Item(1525, 1).Value2 = 4
$ws.Cells.Item(1525, 2).Value2 = 7
$ws.Cells.Item(1525, 3).Value2 = 3
$ws.Cells.Item(1525, 4).Value2 = 13

$ws.Cells.Item(1526, 1).Value2 = 4
$ws.Cells.Item(1526, 2).Value2 = 18
$ws.Cells.Item(1526, 3).Value2 = 1
$ws.Cells.Item(1526, 4).Value2 = 2

$ws.Cells.Item(1527, 1).Value2 = 4
$ws.Cells.Item(1527, 2).Value2 = 3
$ws.Cells.Item(1527, 3).Value2 = 3
$ws.Cells.Item(1527, 4).Value2 = 17

$ws.Cells.Item(1528, 1).Value2 = 9
$ws.Cells.Item(1528, 2).Value2 = 13
$ws.Cells.Item(1528, 3).Value2 = 6
$ws.Cells.Item(1528, 4).Value2 = 7

$ws.Cells.Item(1529, 1).Value2 = 6
$ws.Cells.Item(1529, 2).Value2 = 6
$ws.Cells.Item(1529, 3).Value2 = 5
$ws.Cells.Item(1529, 4).Value2 = 14

$ws.Cells.Item(1530, 1).Value2 = 7
$ws.Cells.Item(1530, 2).Value2 = 14
$ws.Cells.Item(1530, 3).Value2 = 9
$ws.Cells.Item(1530, 4).Value2 = 6

$ws.Cells.Item(1531, 1).Value2 = 6
$ws.Cells.Item(1531, 2).Value2 = 19
$ws.Cells.Item(1531, 3).Value2 = 5
$ws.Cells.Item(1531, 4).Value2 = 1

$ws.Cells.Item(1532, 1).Value2 = 7
$ws.Cells.Item(1532, 2).Value2 = 12
$ws.Cells.Item(1532, 3).Value2 = 5
$ws.Cells.Item(1532, 4).Value2 = 8

$ws.Cells.Item(1533, 1).Value2 = 3
$ws.Cells.Item(1533, 2).Value2 = 18
$ws.Cells.Item(1533, 3).Value2 = 4
$ws.Cells.Item(1533, 4).Value2 = 2

$ws.Cells.Item(1534, 1).Value2 = 3
$ws.Cells.Item(1534, 2).Value2 = 7
$ws.Cells.Item(1534, 3).Value2 = 5
$ws.Cells.Item(1534, 4).Value2 = 13

$ws.Cells.Item(1535, 1).Value2 = 1
$ws.Cells.Item(1535, 2).Value2 = 8
$ws.Cells.Item(1535, 3).Value2 = 3
$ws.Cells.Item(1535, 4).Value2 = 12

$ws.Cells.Item(1536, 1).Value2 = 2
$ws.Cells.Item(1536, 2).Value2 = 5
$ws.Cells.Item(1536, 3).Value2 = 3
$ws.Cells.Item(1536, 4).Value2 = 15

$ws.Cells.Item(1537, 1).Value2 = 4
$ws.Cells.Item(1537, 2).Value2 = 6
$ws.Cells.Item(1537, 3).Value2 = 5
$ws.Cells.Item(1537, 4).Value2 = 14

$ws.Cells.Item(1538, 1).Value2 = 5
$ws.Cells.Item(1538, 2).Value2 = 7
$ws.Cells.Item(1538, 3).Value2 = 6
$ws.Cells.Item(1538, 4).Value2 = 13

$ws.Cells.Item(1539, 1).Value2 = 7
$ws.Cells.Item(1539, 2).Value2 = 14
$ws.Cells.Item(1539, 3).Value2 = 4
$ws.Cells.Item(1539, 4).Value2 = 6

$ws.Activate()
$ws.Range("A1540").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1520
$excel.ActiveWindow.ScrollColumn = 1
